# Add a new "id_posyandu" column to Table1 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one column (Table1: A1:D31 -> A1:E31).
$tbl = $ws.ListObjects.Item(1)
$newCol = $tbl.ListColumns.Add()

# Header + first data row for the new column.
$ws.Cells.Item(1, 5).Value = "id_posyandu"
$ws.Cells.Item(2, 5).Value = 1

# Selection moves to H8 as in the authored workbook.
$ws.Range("H8").Select()
